# Apply a cyclic permutation of the weekly price records in rows 2-7.
# Columns A,B,C,E,F,G,H,I,J,K stay identical (same market/product info).
# Columns D,L,M,N,O,P,Q,R,S,T are rewritten with values "rotated" across rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ D = 44334; L = "Primera"; M = 120; N = 12000; O = 13000; P = 12500; Q = "`$/caja 12 kilos empedrada"; R = "Región de O'Higgins"; S = 1042; T = 12 }
    3 = @{ D = 44698; L = "Primera"; M = 120; N = 16000; O = 17000; P = 16500; Q = "`$/caja 18 kilos granel";     R = "Región de O'Higgins"; S = 917;  T = 18 }
    4 = @{ D = 44344; L = "Primera"; M = 120; N = 13000; O = 14000; P = 13500; Q = "`$/caja 18 kilos granel";     R = "Provincia de Curicó"; S = 750;  T = 18 }
    5 = @{ D = 44316; L = "Primera"; M = 60;  N = 17500; O = 18000; P = 17750; Q = "`$/caja 16 kilos granel";     R = "Región de O'Higgins"; S = 1109; T = 16 }
    6 = @{ D = 44316; L = "Segunda"; M = 40;  N = 16000; O = 16000; P = 16000; Q = "`$/caja 16 kilos granel";     R = "Región de O'Higgins"; S = 1000; T = 16 }
    7 = @{ D = 44330; L = "Primera"; M = 60;  N = 15000; O = 16000; P = 15500; Q = "`$/caja 18 kilos granel";     R = "Provincia de Curicó"; S = 861;  T = 18 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value  = $vals.D   # D: Fecha
    $ws.Cells.Item($row, 12).Value = $vals.L   # L: Calidad
    $ws.Cells.Item($row, 13).Value = $vals.M   # M: Volumen
    $ws.Cells.Item($row, 14).Value = $vals.N   # N: Precio minimo
    $ws.Cells.Item($row, 15).Value = $vals.O   # O: Precio maximo
    $ws.Cells.Item($row, 16).Value = $vals.P   # P: Precio promedio ponderado
    $ws.Cells.Item($row, 17).Value = $vals.Q   # Q: Unidad de comercializacion
    $ws.Cells.Item($row, 18).Value = $vals.R   # R: Origen
    $ws.Cells.Item($row, 19).Value = $vals.S   # S: Precio $/Kg
    $ws.Cells.Item($row, 20).Value = $vals.T   # T: Kg / unidad
}
